# Update Data by bot, scripted by HH
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 must stay text ("001", not the number 1) even though it looks numeric.
# Force text formatting while the value is assigned, then drop the
# NumberFormat override again so the cell's style is left untouched.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"
$ws.Range("J2").ClearFormats()

$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 74617232.42
$ws.Range("P2").Value = 840125824.51
$ws.Range("Q2").Value = 658192468.78
$ws.Range("R2").Value = -40.5801766401
$ws.Range("S2").Value = 551763159.7
$ws.Range("T2").Value = 551763159.7
$ws.Range("U2").Value = -43.3845619894
$ws.Range("V2").Value = 33726088.97
$ws.Range("W2").Value = 97328288.76000001
$ws.Range("X2").Value = -66147911.6
$ws.Range("Y2").Value = 31827798.96
$ws.Range("Z2").Value = 19131818.19
$ws.Range("AA2").Value = -51874692.95

$ws.Range("AG2").Value = 13512314.88

$ws.Range("AP2").Value = -36.5598255413
$ws.Range("AQ2").Value = 154.399145628435
$ws.Range("AR2").Value = 182.58
$ws.Range("AS2").Value = 73059847.69
$ws.Range("AT2").Value = 183.498554965292
